# Fruta / hortaliza, semanal
# Insert a new data row above row 242 (pushing existing rows 242-263 down
# to 243-264) and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 242..263 down to 243..264, carrying formatting along.
$ws.Rows.Item(242).EntireRow.Insert()

# Fill in the brand-new row 242 with this week's price entry.
$ws.Cells.Item(242, 1).Value  = 7
$ws.Cells.Item(242, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(242, 3).Value  = "Ñuble"
$ws.Cells.Item(242, 4).Value  = 45132
$ws.Cells.Item(242, 5).Value  = 16
$ws.Cells.Item(242, 6).Value  = "Fruta"
$ws.Cells.Item(242, 7).Value  = 100109
$ws.Cells.Item(242, 8).Value  = "Uva"
$ws.Cells.Item(242, 9).Value  = 100109001
$ws.Cells.Item(242, 10).Value = "Uva"
$ws.Cells.Item(242, 11).Value = "Crimpson Seedless"
$ws.Cells.Item(242, 12).Value = "Primera"
$ws.Cells.Item(242, 13).Value = 60
$ws.Cells.Item(242, 14).Value = 12000
$ws.Cells.Item(242, 15).Value = 12000
$ws.Cells.Item(242, 16).Value = 12000
$ws.Cells.Item(242, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(242, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(242, 19).Value = 1500
$ws.Cells.Item(242, 20).Value = 8

# Make sure the new date cell uses the same date style as its neighbours.
$ws.Cells.Item(242, 4).NumberFormat = $ws.Cells.Item(243, 4).NumberFormat
